$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New A/B column values (rows 1-32)
$aVals = @(-0.11550149516097719, -0.093335444317927063, -0.075334803093461744, -0.067169995774468916, -0.063668358716949314, -0.041402982331161908, -0.030987864095696427, -0.020892209160482622, -0.018737361410777975, -0.016615416177847564, -0.013607792428008914, -0.010096387482867453, -0.0065212121333653172, -0.0090817727691954175, -0.0080530463533907692, -0.0060344892452137699, -0.0040035356932675015, -0.016104988816852739, -0.012091530277658968, -0.008016810508374661, -0.0040056818418170437, -0.045703650338701607, -0.040492785238618545, -0.020097817630094106, -0.018554195553683428, -0.016006032665087488, -0.013447498339905994, -0.011122255949429061, -0.081347710440050847, -0.021166827626636131, -0.014023023405266954, -0.0040013085635077772)
$bVals = @(0.11544043310843932, 0.093164298046992222, 0.075169995734574613, 0.066668358696714947, 0.061957177864607083, 0.040987864040467059, 0.030892209103662527, 0.020737361386611752, 0.01861541615217055, 0.016607792398044552, 0.013596387450410141, 0.010021212099727173, 0.0064931800795795525, 0.0090530463298028607, 0.0080344892168735527, 0.0060035356643011184, 0.0039999999622590821, 0.016091530259814135, 0.012016810489242413, 0.0080056818224853998, 0.0039999999805155895, 0.04549278521200506, 0.040097817536507208, 0.019999999905198074, 0.018506032642488179, 0.015947498316982944, 0.013122255927460635, 0.010918640856397488, 0.081166827351850657, 0.021023023356466553, 0.014001308501502052, 0.0039999999636890493)

for ($i = 0; $i -lt 32; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}

# Updated column widths (nearest value Excel's pixel-quantized ColumnWidth
# property can represent for the target stored widths of 15.42578125 / 14.7109375)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334
